$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 45946): B2, C2, F2 currently empty, now get values
$ws.Range("B2").Value = 5919.77118360843
$ws.Range("C2").Value = 6143.0263933066
$ws.Range("F2").Value = 75.4315641124237

# Row 3 (A3 = 45947)
$ws.Range("B3").Value = 6034.6354574457
$ws.Range("C3").Value = 5555.29716325231
$ws.Range("F3").Value = 200.089205533609

# Row 4 (A4 = 45948)
$ws.Range("B4").Value = 2015.18486681053
$ws.Range("C4").Value = 4082.46831061959
$ws.Range("F4").Value = 126.424967408711

# Row 5 (A5 = 45949)
$ws.Range("B5").Value = 2022.9840958268
$ws.Range("C5").Value = 4206.33690883181
$ws.Range("F5").Value = 133.418827708542

# Row 6 (A6 = 45950)
$ws.Range("B6").Value = 7124.77498284078
$ws.Range("C6").Value = 6948.93019712549
$ws.Range("F6").Value = 289.710888095196

# Row 7 (A7 = 45951)
$ws.Range("B7").Value = 7405.98294146144
$ws.Range("C7").Value = 6992.29409233059
$ws.Range("F7").Value = 245.618293786214
